$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pareto exponent extrapolation values were incorrectly rounded (left blank).
# Fill in the correct computed values for rows 12, 13, 15, 18 (columns I-L).

$ws.Range("I12").Value = -0.009480606601132577
$ws.Range("J12").Value = 0.01578916715912978
$ws.Range("K12").Value = -0.2836171216747701
$ws.Range("L12").Value = 2.311125364129039

$ws.Range("I13").Value = 0.008713973304380595
$ws.Range("J13").Value = 0.01551982827831129
$ws.Range("K13").Value = -0.5380701125052759
$ws.Range("L13").Value = 2.392278948236461

$ws.Range("I15").Value = 0.06622089617124935
$ws.Range("J15").Value = 0.03672917298683431
$ws.Range("K15").Value = 0.6725847961907231
$ws.Range("L15").Value = 3.168892690376738

$ws.Range("I18").Value = 0.1697846028918021
$ws.Range("J18").Value = 0.07334250380262709
$ws.Range("K18").Value = 0.4683031253288178
$ws.Range("L18").Value = 2.011110609254798
